$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the Price/Volume columns so numeric-looking
# strings (e.g. "1.00", "66.655.87") are preserved verbatim as text,
# matching the source data which stores these as inline strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "66.655.87"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "2.521.87"
$ws.Range("E3").Value = "  -4.06%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "584.62"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").Value = "171.67"
$ws.Range("E6").Value = "  +2.11%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("D9").Value = "2.521.68"
$ws.Range("E9").Value = "  -4.10%  "
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "0.350"
$ws.Range("E12").Value = "  -4.03%  "
$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D13").Value = "5.13"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").Value = "26.74"
$ws.Range("E14").Value = "  -3.28%  "
$ws.Range("D15").Value = "2.980.48"
$ws.Range("E15").Value = "  -4.40%  "
$ws.Range("E16").Value = "  -2.99%  "
$ws.Range("D17").Value = "66.581.80"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "2.503.28"
$ws.Range("E18").Value = "  -4.57%  "
$ws.Range("D19").Value = "7.85"
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").Value = "11.33"
$ws.Range("E20").Value = "  -5.81%  "
$ws.Range("D21").Value = "348.03"
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("D23").Value = "4.65"
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  +2.19%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "70.21"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "9.96"
$ws.Range("E27").Value = "  -3.91%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "2.634.48"
$ws.Range("E29").Value = "  -4.72%  "
$ws.Range("D30").Value = "0.0₃0979"
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("D31").Value = "526.60"
$ws.Range("E31").Value = "  -3.89%  "
$ws.Range("D32").Value = "8.13"
$ws.Range("E32").Value = "  +2.12%  "
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("E34").Value = "  -2.96%  "
$ws.Range("E35").Value = "  -4.62%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -2.77%  "
$ws.Range("D38").Value = "157.58"
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("D39").Value = "18.65"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").Value = "18.40"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("E42").Value = "  -1.08%  "
$ws.Range("D43").Value = "5.10"
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").Value = "  +2.89%  "
$ws.Range("D46").Value = "39.49"
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("D47").Value = "149.52"
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "0.0₆0270"
$ws.Range("E51").Value = "  -10.61%  "
